$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = "Datos actualizados a 26 de Marzo de 2020 a las 09:42"

# Row 26
$ws.Range("F26").Value = 34

# Row 34
$ws.Range("B34").Value = 1085
$ws.Range("C34").Value = 34
$ws.Range("E34").Value = 1063
$ws.Range("G34").Value = 1
$ws.Range("H34").Value = 15

# Row 35
$ws.Range("D35").Value = 88
$ws.Range("E35").Value = 953
$ws.Range("F35").Value = 4

# Row 40
$ws.Range("E40").Value = 762
$ws.Range("G40").Value = 1
$ws.Range("H40").Value = 23

# Row 45
$ws.Range("A45").Value = "Filipinas"
$ws.Range("B45").Value = 707
$ws.Range("C45").Value = 71
$ws.Range("D45").Value = 28
$ws.Range("E45").Value = 634
$ws.Range("F45").Value = 1
$ws.Range("G45").Value = 7
$ws.Range("H45").Value = 45

# Row 46
$ws.Range("A46").Value = "India"
$ws.Range("B46").Value = 681
$ws.Range("C46").Value = 24
$ws.Range("D46").Value = 43
$ws.Range("E46").Value = 625
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 1
$ws.Range("H46").Value = 13

# Row 52
$ws.Range("F52").Value = 14

# Row 71
$ws.Range("A71").Value = "Letonia"
$ws.Range("B71").Value = 244
$ws.Range("C71").Value = 23
$ws.Range("D71").Value = 1
$ws.Range("E71").Value = 243
$ws.Range("F71").Value = 0
$ws.Range("H71").Value = 0

# Row 72
$ws.Range("A72").Value = "Bulgaria"
$ws.Range("B72").Value = 243
$ws.Range("C72").Value = 1
$ws.Range("D72").Value = 4
$ws.Range("E72").Value = 236
$ws.Range("F72").Value = 8
$ws.Range("H72").Value = 3

# Row 73
$ws.Range("A73").Value = "Marruecos"
$ws.Range("B73").Value = 225
$ws.Range("D73").Value = 7
$ws.Range("E73").Value = 212
$ws.Range("F73").Value = 1
$ws.Range("H73").Value = 6

# Row 77
$ws.Range("A77").Value = "Kuwait"
$ws.Range("B77").Value = 208
$ws.Range("C77").Value = 13
$ws.Range("D77").Value = 49
$ws.Range("E77").Value = 159
$ws.Range("F77").Value = 7
$ws.Range("H77").Value = 0

# Row 78
$ws.Range("A78").Value = "Costa Rica"
$ws.Range("B78").Value = 201
$ws.Range("D78").Value = 2
$ws.Range("E78").Value = 197
$ws.Range("F78").Value = 4
$ws.Range("H78").Value = 2

# Row 89
$ws.Range("F89").Value = 0

# Row 117
$ws.Range("A117").Value = "Consejo Danes para los Refugiados"
$ws.Range("C117").Value = 3
$ws.Range("D117").Value = 0
$ws.Range("G117").Value = 1
$ws.Range("H117").Value = 3

# Row 118
$ws.Range("A118").Value = "Nigeria"
$ws.Range("B118").Value = 51
$ws.Range("D118").Value = 2
$ws.Range("E118").Value = 48
$ws.Range("H118").Value = 1
